# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets.
# On "展览" the updated rows are 4,6,7,8,10,11,17,19,20,26,27,28,29,30,33,38,39.
# On "全部类型" the same rows apply except the last four are shifted down by
# two (32,35,40,41) because that sheet has two extra data rows.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览" = @{
        4  = 144
        6  = 18159
        7  = 364
        8  = 261
        10 = 6845
        11 = 688
        17 = 156
        19 = 227
        20 = 56
        26 = 987
        27 = 126
        28 = 5165
        29 = 535
        30 = 35
        33 = 12073
        38 = 3920
        39 = 302
    }
    "全部类型" = @{
        4  = 144
        6  = 18159
        7  = 364
        8  = 261
        10 = 6845
        11 = 688
        17 = 156
        19 = 227
        20 = 56
        26 = 987
        27 = 126
        28 = 5165
        29 = 535
        32 = 35
        35 = 12073
        40 = 3920
        41 = 302
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
